$wb = $excel.ActiveWorkbook

# Add the new "table_definitions" worksheet and move it to be the first tab
$defs = $wb.Worksheets.Add()
$defs.Name = "table_definitions"
$defs.Move($wb.Worksheets.Item(1))

# Header row (row 1)
$defs.Range("A1").Value = "mapping_file_name"
$defs.Range("B1").Value = "entity_name"
$defs.Range("C1").Value = "required_entities"
$defs.Range("D1").Value = "destination_table_name"
$defs.Range("E1").Value = "table_type"
$defs.Range("F1").Value = "source_table_name"
$defs.Range("G1").Value = "casrec_conditions"
$defs.Range("H1").Value = "source_table_additional_columns"

# Data row (row 2)
$defs.Range("A2").Value = "crec_persons"
$defs.Range("B2").Value = "crec"
$defs.Range("C2").Value = "client"
$defs.Range("D2").Value = "persons"
$defs.Range("E2").Value = "data"
$defs.Range("F2").Value = "crec"

# Explicit font formatting on a subset of the header cells
$defs.Range("A1,D1,E1,F1,G1").Font.Name = "Arial"
$defs.Range("A1,D1,E1,F1,G1").Font.Size = 10

# Mirror the author's final on-sheet selection (whole rows 1-2 selected)
$defs.Rows("1:2").Select()
